$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text
$ws.Range("B1").Value = "Value (g)"

# Rows whose numeric values become text strings holding the same displayed
# value (rows 11, 23, 29, 31 are left untouched - still plain numbers).
$textValues = [ordered]@{
    2  = "26.9"
    3  = "40.8"
    4  = "209.9"
    5  = "45.4"
    6  = "105.8"
    7  = "31.7"
    8  = "54.3"
    9  = "44.2"
    10 = "36.9"
    12 = "75.2"
    13 = "87.8"
    14 = "44.1"
    15 = "44.6"
    16 = "33.7"
    17 = "57.5"
    18 = "166.3"
    19 = "77.8"
    20 = "27.5"
    21 = "81.2"
    22 = "31.3"
    24 = "27.5"
    25 = "69.9"
    26 = "44.2"
    27 = "7.8"
    28 = "118.8"
    30 = "21.6"
    32 = "97.7"
    33 = "46.2"
    34 = "97.6"
    35 = "77.7"
}

foreach ($row in $textValues.Keys) {
    $cell = $ws.Cells.Item($row, 2)
    # Write a text-producing formula, then paste-special as values so the
    # resulting cell holds a literal text value (no style/number-format
    # change, unlike assigning a numeric-looking string to .Value).
    $cell.Formula = '="' + $textValues[$row] + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = 0

# Autofit column A width
$ws.Columns.Item(1).AutoFit() | Out-Null

# Set active selection to B2
$ws.Range("B2").Select() | Out-Null
